$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 6349.9375
$ws.Range("I28").Value = 843.8889
$ws.Range("J28").Value = 13429.143
$ws.Range("K28").Value = 843.8889
$ws.Range("L28").Value = 13429.143
$ws.Range("M28").Value = -358.8889
$ws.Range("N28").Value = -14399.143
# Row 33
$ws.Range("H33").Value = 314.375
$ws.Range("I33").Value = 314.375
$ws.Range("K33").Value = 314.375
$ws.Range("M33").Value = -85.375
# Row 105
$ws.Range("H105").Value = 34335.5
$ws.Range("J105").Value = 34335.5
$ws.Range("L105").Value = 34335.5
$ws.Range("N105").Value = -41323.5
# Row 107
$ws.Range("H107").Value = 175.5
$ws.Range("I107").Value = 175.5
$ws.Range("K107").Value = 175.5
$ws.Range("M107").Value = 1744.5
# Row 132
$ws.Range("H132").Value = 6670.1
$ws.Range("I132").Value = 2185.8572
$ws.Range("K132").Value = 6557.571599999999
$ws.Range("M132").Value = -4027.571599999999
# Row 137
$ws.Range("H137").Value = 2914.3809
$ws.Range("I137").Value = 1688.8889
$ws.Range("J137").Value = 3833.5
$ws.Range("K137").Value = 5066.6667
$ws.Range("L137").Value = 11500.5
$ws.Range("M137").Value = -2516.6667
$ws.Range("N137").Value = -16600.5
# Row 138
$ws.Range("H138").Value = 2475.25
$ws.Range("J138").Value = 3490
$ws.Range("L138").Value = 10470
$ws.Range("N138").Value = -20750

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 500
$ws.Range("J10").Value = 500
$ws.Range("L10").Value = 500
$ws.Range("N10").Value = -840
# Row 13
$ws.Range("H13").Value = 1549
$ws.Range("I13").Value = 1349
$ws.Range("J13").Value = 1649
$ws.Range("K13").Value = 1349
$ws.Range("L13").Value = 1649
$ws.Range("M13").Value = -1205
$ws.Range("N13").Value = -1937
# Row 61
$ws.Range("H61").Value = 4299
$ws.Range("I61").Value = 3833.3333
$ws.Range("J61").Value = 4997.5
$ws.Range("K61").Value = 3833.3333
$ws.Range("L61").Value = 4997.5
$ws.Range("M61").Value = -3621.3333
$ws.Range("N61").Value = -5421.5
# Row 74
$ws.Range("I74").Value = 9999
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 9999
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -9125
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("I77").Value = 9999
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 49995
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -45627
$ws.Range("N77").ClearContents()
# Row 132
$ws.Range("H132").Value = 2291.8572
$ws.Range("I132").Value = 2473.8333
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 7421.499899999999
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -4891.499899999999
$ws.Range("N132").Value = -8660
# Row 136
$ws.Range("H136").Value = 4299
$ws.Range("I136").Value = 3833.3333
$ws.Range("J136").Value = 4997.5
$ws.Range("K136").Value = 11499.9999
$ws.Range("L136").Value = 14992.5
$ws.Range("M136").Value = -8949.999899999999
$ws.Range("N136").Value = -20092.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 393.75
$ws.Range("I5").Value = 191.66667
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 191.66667
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -78.66667000000001
$ws.Range("N5").Value = -1226
# Row 94
$ws.Range("H94").Value = 335.36365
$ws.Range("I94").Value = 268.9
$ws.Range("K94").Value = 268.9
$ws.Range("M94").Value = 182.1
# Row 99
$ws.Range("H99").Value = 1398.2
$ws.Range("I99").Value = 1245
$ws.Range("J99").Value = 2011
$ws.Range("K99").Value = 1245
$ws.Range("L99").Value = 2011
$ws.Range("M99").Value = 253
$ws.Range("N99").Value = -5007

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 67.40000000000001
$ws.Range("I7").Value = 70.28570999999999
$ws.Range("K7").Value = 70.28570999999999
$ws.Range("M7").Value = 42.71429000000001
# Row 22
$ws.Range("H22").Value = 1472.5
$ws.Range("I22").Value = 826
$ws.Range("J22").Value = 2334.5
$ws.Range("K22").Value = 826
$ws.Range("L22").Value = 2334.5
$ws.Range("M22").Value = -476
$ws.Range("N22").Value = -3034.5
# Row 31
$ws.Range("H31").Value = 7338.1665
$ws.Range("I31").Value = 3624
$ws.Range("J31").Value = 8688.772000000001
$ws.Range("K31").Value = 3624
$ws.Range("L31").Value = 8688.772000000001
$ws.Range("M31").Value = -3329
$ws.Range("N31").Value = -9278.772000000001
# Row 32
$ws.Range("H32").Value = 1371.25
$ws.Range("I32").Value = 1138.5714
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 1138.5714
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -822.5714
$ws.Range("N32").Value = -3632
# Row 34
$ws.Range("H34").Value = 7338.1665
$ws.Range("I34").Value = 3624
$ws.Range("J34").Value = 8688.772000000001
$ws.Range("K34").Value = 3624
$ws.Range("L34").Value = 8688.772000000001
$ws.Range("M34").Value = -3422
$ws.Range("N34").Value = -9092.772000000001
# Row 41
$ws.Range("H41").Value = 64596
$ws.Range("J41").Value = 74364.5
$ws.Range("L41").Value = 74364.5
$ws.Range("N41").Value = -75220.5
# Row 58
$ws.Range("H58").Value = 4885.091
$ws.Range("I58").Value = 1402.4
$ws.Range("K58").Value = 1402.4
$ws.Range("M58").Value = -1199.4
# Row 99
$ws.Range("H99").Value = 4166.1665
$ws.Range("J99").Value = 4166.1665
$ws.Range("L99").Value = 4166.1665
$ws.Range("N99").Value = -7162.1665
# Row 126
$ws.Range("H126").Value = 4166.1665
$ws.Range("J126").Value = 4166.1665
$ws.Range("L126").Value = 12498.4995
$ws.Range("N126").Value = -17438.4995
# Row 132
$ws.Range("H132").Value = 4998.8
$ws.Range("I132").Value = 4249.1665
$ws.Range("K132").Value = 12747.4995
$ws.Range("M132").Value = -10217.4995
# Row 136
$ws.Range("H136").Value = 4885.091
$ws.Range("I136").Value = 1402.4
$ws.Range("K136").Value = 4207.200000000001
$ws.Range("M136").Value = -1657.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 260.2353
$ws.Range("J12").Value = 274.5
$ws.Range("L12").Value = 823.5
$ws.Range("N12").Value = -1169.5
# Row 107
$ws.Range("H107").Value = 749.1667
$ws.Range("I107").Value = 497.5
$ws.Range("K107").Value = 1492.5
$ws.Range("M107").Value = 427.5
# Row 124
$ws.Range("H124").Value = 1399.5
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 125
$ws.Range("H125").Value = 7998.5
$ws.Range("I125").Value = 7998.5
$ws.Range("K125").Value = 23995.5
$ws.Range("M125").Value = -19075.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 17600526
$ws.Range("I11").Value = 17647544
$ws.Range("K11").Value = 17647544
$ws.Range("M11").Value = -17647405
# Row 102
$ws.Range("H102").Value = 1497.3077
$ws.Range("I102").Value = 1588.75
$ws.Range("J102").Value = 400
$ws.Range("K102").Value = 1588.75
$ws.Range("L102").Value = 400
$ws.Range("M102").Value = 33.25
$ws.Range("N102").Value = -3644
# Row 122
$ws.Range("H122").Value = 2625.1
$ws.Range("I122").Value = 1126.3334
$ws.Range("J122").Value = 4873.25
$ws.Range("K122").Value = 3379.0002
$ws.Range("L122").Value = 14619.75
$ws.Range("M122").Value = -929.0001999999999
$ws.Range("N122").Value = -19519.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 18000
$ws.Range("I14").Value = 18000
$ws.Range("K14").Value = 18000
$ws.Range("M14").Value = -17828
# Row 22
$ws.Range("H22").Value = 878
$ws.Range("I22").Value = 650
$ws.Range("J22").Value = 935
$ws.Range("K22").Value = 650
$ws.Range("L22").Value = 935
$ws.Range("M22").Value = -355
$ws.Range("N22").Value = -1525
# Row 27
$ws.Range("H27").Value = 878
$ws.Range("I27").Value = 650
$ws.Range("J27").Value = 935
$ws.Range("K27").Value = 650
$ws.Range("L27").Value = 935
$ws.Range("M27").Value = -543
$ws.Range("N27").Value = -1149
# Row 46
$ws.Range("H46").Value = 5431.35
$ws.Range("I46").Value = 3847.4443
$ws.Range("J46").Value = 6727.273
$ws.Range("K46").Value = 3847.4443
$ws.Range("L46").Value = 6727.273
$ws.Range("M46").Value = -3659.4443
$ws.Range("N46").Value = -7103.273
# Row 93
$ws.Range("H93").Value = 2201
$ws.Range("I93").Value = 2499
$ws.Range("J93").Value = 2052
$ws.Range("K93").Value = 2499
$ws.Range("L93").Value = 2052
$ws.Range("M93").Value = -1251
$ws.Range("N93").Value = -4548
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 132
$ws.Range("H132").Value = 4409
$ws.Range("I132").Value = 3476.3333
$ws.Range("K132").Value = 10428.9999
$ws.Range("M132").Value = -7898.999899999999
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 320.6757
$ws.Range("I14").Value = 328.26086
$ws.Range("J14").Value = 308.2143
$ws.Range("K14").Value = 328.26086
$ws.Range("L14").Value = 308.2143
$ws.Range("M14").Value = -160.26086
$ws.Range("N14").Value = -644.2143
# Row 126
$ws.Range("H126").Value = 5140
$ws.Range("I126").Value = 1992.4
$ws.Range("K126").Value = 5977.200000000001
$ws.Range("M126").Value = -3507.200000000001
# Row 136
$ws.Range("H136").Value = 3599
$ws.Range("I136").Value = 4898.3335
$ws.Range("J136").Value = 1650
$ws.Range("K136").Value = 14695.0005
$ws.Range("L136").Value = 4950
$ws.Range("M136").Value = -12145.0005
$ws.Range("N136").Value = -10050
